$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("idef")

# Restrict the excludeDisaggregation value for the "tasa_ocupacion" row (row 2)
# to only the pair indicated: "ethnicity,disability,migrante,area"
$ws.Range("F2").Value = "ethnicity,disability,migrante,area"

# Update the active selection to match the edited workbook state
$ws.Range("D7").Select()
